$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

# Title heading and the duplicated bold "Play Combat Masters..." line near the end
Replace-Text "Play Combat Masters Free - Slot Game Review" "Play Combat Masters Free Game - Review & Ratings"

# "What we like" bullets
Replace-Text "Engaging gameplay system with cascade fall of symbols" "Engaging gameplay with cascade fall of symbols"
Replace-Text "8 different levels adding to the excitement" "Eight different levels provide variety"
Replace-Text "Wild Power feature with extra Wilds adds value to the gameplay" "Gameplay-enhancing functions and Wild Power feature"
Replace-Text "Modern graphics with a fascinating theme and environment" "Modern graphics and captivating design"

# "What we don't like" bullets
Replace-Text "Significant awards might take some time due to medium volatility" "Significant awards may take time"
Replace-Text "Fixed coin value limits flexibility in betting" "Medium volatility may not be suitable for all players"

# Meta description italic text
Replace-Text "Experience the engaging gameplay system of Combat Masters! Read our review and play for free to enjoy medium volatility and fascinating graphics." "Play Combat Masters for free and enjoy engaging gameplay, modern graphics, and captivating design."
